$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy style (bold font + border + center/top alignment) from A16 down to new label cells A17:A19
$ws.Range("A16").Copy()
$ws.Range("A17:A19").PasteSpecial(-4122)

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("C10").Value = 1.067466239315418
$ws.Range("D10").Value = 0.8773520961605883
$ws.Range("E10").Value = 1.070699185305034
$ws.Range("F10").Value = 0.9439715217766051
$ws.Range("G10").Value = 1.067466239315418
$ws.Range("H10").Value = 0.8773520961605883
$ws.Range("I10").Value = 1.055506493621435
$ws.Range("J10").Value = 0.9784606975843092
$ws.Range("K10").Value = 1.003181088641286
$ws.Range("L10").Value = 0.9002881726143893
$ws.Range("M10").Value = 1.067466239315418
$ws.Range("N10").Value = 0.974025640732811
$ws.Range("O10").Value = 0.9898722606394112
$ws.Range("P10").Value = 0.9871156868773832

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("C11").Value = 1.077804198382943
$ws.Range("D11").Value = 0.3440270037562719
$ws.Range("E11").Value = 1.403425352726551
$ws.Range("F11").Value = 0.8376311809785424
$ws.Range("G11").Value = 1.077804198382943
$ws.Range("H11").Value = 0.3440270037562719
$ws.Range("I11").Value = 1.291076887109758
$ws.Range("J11").Value = 0.9703932391322229
$ws.Range("K11").Value = 1.043558715844727
$ws.Range("L11").Value = 0.573913743313366
$ws.Range("M11").Value = 1.077804198382943
$ws.Range("N11").Value = 0.8737261782414114
$ws.Range("O11").Value = 0.9157219339610772
$ws.Range("P11").Value = 0.9427287901555478

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("C12").Value = 1.078393670575991
$ws.Range("D12").Value = 0.3442814100132307
$ws.Range("E12").Value = 1.402481029769529
$ws.Range("F12").Value = 0.8378528482072
$ws.Range("G12").Value = 1.078393670575991
$ws.Range("H12").Value = 0.3442814100132307
$ws.Range("I12").Value = 1.290811732794739
$ws.Range("J12").Value = 0.9701517281235125
$ws.Range("K12").Value = 1.043881557406747
$ws.Range("L12").Value = 0.5739688882246932
$ws.Range("M12").Value = 1.078393670575991
$ws.Range("N12").Value = 0.87338121989138
$ws.Range("O12").Value = 0.9157522396414877
$ws.Range("P12").Value = 0.9427278581394551

$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("C13").Value = 1.077482013124339
$ws.Range("D13").Value = 0.3439787073714256
$ws.Range("E13").Value = 1.403790993823035
$ws.Range("F13").Value = 0.8375590272116238
$ws.Range("G13").Value = 1.077482013124339
$ws.Range("H13").Value = 0.3439787073714256
$ws.Range("I13").Value = 1.29123217955278
$ws.Range("J13").Value = 0.9705260111528964
$ws.Range("K13").Value = 1.043381397333751
$ws.Range("L13").Value = 0.5737050102409803
$ws.Range("M13").Value = 1.077482013124339
$ws.Range("N13").Value = 0.8738848505972303
$ws.Range("O13").Value = 0.915702685382606
$ws.Range("P13").Value = 0.942706917476354

$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("C14").Value = 1.309716000000002
$ws.Range("D14").Value = 0.3944600000000001
$ws.Range("E14").Value = 1.108312000000001
$ws.Range("F14").Value = 0.9071920000000002
$ws.Range("G14").Value = 1.309716000000002
$ws.Range("H14").Value = 0.3944600000000001
$ws.Range("I14").Value = 1.180672
$ws.Range("J14").Value = 0.8754920000000015
$ws.Range("K14").Value = 1.162468000000001
$ws.Range("L14").Value = 0.6492720000000026
$ws.Range("M14").Value = 1.309716000000002
$ws.Range("N14").Value = 0.7513860000000003
$ws.Range("O14").Value = 0.9299200000000007
$ws.Range("P14").Value = 0.9484480000000011

$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("C15").Value = 1.524012500000001
$ws.Range("D15").Value = 0.44
$ws.Range("E15").Value = 0.84
$ws.Range("F15").Value = 0.9701375000000003
$ws.Range("G15").Value = 1.524012500000001
$ws.Range("H15").Value = 0.44
$ws.Range("I15").Value = 1.079724999999999
$ws.Range("J15").Value = 0.7901375000000004
$ws.Range("K15").Value = 1.271237499999998
$ws.Range("L15").Value = 0.72
$ws.Range("M15").Value = 1.524012500000001
$ws.Range("N15").Value = 0.64
$ws.Range("O15").Value = 0.9435375000000004
$ws.Range("P15").Value = 0.9544062499999999

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "Rotation-60detTilt"
$ws.Range("C16").Value = 1.299947330559997
$ws.Range("D16").Value = 0.6708016477184017
$ws.Range("E16").Value = 0.907142698598401
$ws.Range("F16").Value = 0.9796383726591984
$ws.Range("G16").Value = 1.299947330559997
$ws.Range("H16").Value = 0.6708016477184017
$ws.Range("I16").Value = 1.045323299123201
$ws.Range("J16").Value = 0.8776181765120009
$ws.Range("K16").Value = 1.153602846003199
$ws.Range("L16").Value = 0.8333974270975991
$ws.Range("M16").Value = 1.299938626150397
$ws.Range("N16").Value = 0.7889721731584014
$ws.Range("O16").Value = 0.9643825123839993
$ws.Range("P16").Value = 0.9709339747839997

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C17").Value = 0.9932715404578597
$ws.Range("D17").Value = 0.9938945841856551
$ws.Range("E17").Value = 0.9956860940398954
$ws.Range("F17").Value = 0.9934079569945276
$ws.Range("G17").Value = 0.9932715404578597
$ws.Range("H17").Value = 0.9938945841856551
$ws.Range("I17").Value = 0.9946744096261005
$ws.Range("J17").Value = 0.9956769371794372
$ws.Range("K17").Value = 0.9941687112364865
$ws.Range("L17").Value = 0.9927591642822557
$ws.Range("M17").Value = 0.9932631326099548
$ws.Range("N17").Value = 0.9947903391127753
$ws.Range("O17").Value = 0.9940650439194845
$ws.Range("P17").Value = 0.9941924247502772

$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("C18").Value = 1.017839736054624
$ws.Range("D18").Value = 1.028569525148503
$ws.Range("E18").Value = 0.95585870048236
$ws.Range("F18").Value = 1.003196286303925
$ws.Range("G18").Value = 1.017839736054624
$ws.Range("H18").Value = 1.028569525148503
$ws.Range("I18").Value = 0.9739405772916858
$ws.Range("J18").Value = 0.9843320482359066
$ws.Range("K18").Value = 1.002172672315638
$ws.Range("L18").Value = 1.01296671287698
$ws.Range("M18").Value = 1.01787297336068
$ws.Range("N18").Value = 0.9922141128154314
$ws.Range("O18").Value = 1.001366061997353
$ws.Range("P18").Value = 0.9973595323387028

$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C19").Value = 0.9486806746275404
$ws.Range("D19").Value = 1.046172334010831
$ws.Range("E19").Value = 0.9766981258125575
$ws.Range("F19").Value = 1.011031250459383
$ws.Range("G19").Value = 0.9486806746275404
$ws.Range("H19").Value = 1.046172334010831
$ws.Range("I19").Value = 0.9693539070546845
$ws.Range("J19").Value = 1.009117548961016
$ws.Range("K19").Value = 0.9819569199129816
$ws.Range("L19").Value = 1.034371278075232
$ws.Range("M19").Value = 0.9486586716556321
$ws.Range("N19").Value = 1.011435229911694
$ws.Range("O19").Value = 0.9956455962275779
$ws.Range("P19").Value = 0.9971727548642781
